$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear columns AB:AK (No. of Sites/bldg breakdown columns) and AM (DIFFERENCE)
# for data rows 2-18, keeping AL (PREVIOUS ACCOMPLISHMENT) intact.
for ($r = 2; $r -le 18; $r++) {
    $ws.Range("AB$r`:AK$r").ClearContents()
    $ws.Range("AM$r").ClearContents()
}
